$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Email" column (L), shifting
# Email (and everything to its right) one column over to M. Excel keeps
# the header style (bold/border) on the new blank L1 cell automatically.
$ws.Columns("L:L").Insert()

# Header for the newly inserted column.
$ws.Range("L1").Value = "Resenha Regime de Metas"

# Per-student "Resenha Regime de Metas" grades (column L). Rows whose
# other grade columns are also blank for this assignment are left empty.
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 7
$ws.Range("L5").Value = 7
$ws.Range("L6").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = 7
$ws.Range("L10").Value = 7
$ws.Range("L11").Value = 0
$ws.Range("L13").Value = 5
$ws.Range("L14").Value = 5
$ws.Range("L16").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("L18").Value = 7
$ws.Range("L19").Value = 10
$ws.Range("L20").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("L22").Value = 7
$ws.Range("L23").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("L31").Value = 7
$ws.Range("L32").Value = 10
$ws.Range("L33").Value = 7
$ws.Range("L34").Value = 10
$ws.Range("L35").Value = 7
$ws.Range("L36").Value = 3
$ws.Range("L37").Value = 0
